# Remove the log entry that was at row 19 (Student ID 201438, logged at
# 11:25:50). Deleting the worksheet row shifts every row below it up by
# one, which reproduces the row-by-row "off by one" shift seen across
# rows 19-113 in the diff, and drops the former last row (113) so the
# sheet ends at row 112.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Delete()
